# Apply "a lot of updates to strain_comparison" changes to the Allele sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Allele")

# ---------------------------------------------------------------------------
# 1. Small text / value edits to existing rows
# ---------------------------------------------------------------------------

# Row 37 (AlleleID 37): me2-mero genotype gains PgapD promoter prefix
$ws.Range("F37").Value = "∆me::PgapD-cat-hpt"

# Row 109 (AlleleID 110): adhE6 plasmid genotype, fix PgapDH -> PgapD
$ws.Range("F109").Value = "pSH007(PgapD-adhE)"

# Row 110 (AlleleID 111): pta1-mero genotype gains PgapD promoter prefix
$ws.Range("F110").Value = "∆pta::PgapD-cat-hpt"

# Row 112 (AlleleID 113): tscEtoh2 -> tscEtoh2int, now an integrated (not
# plasmid-borne) construct
$ws.Range("C112").Value = "tscEtoh2int"
$ws.Range("D112").Value = "ethanol"
$ws.Range("F112").Value = "P2638::pSH64(P2638-nfnAB(Tsc)-adhA(Tsc))"
$ws.Range("G112").Value = $false
$ws.Range("N112").Value = "pSH064 integrated in chromosome at Clo1313_2638 locus"

# Row 121 (AlleleID 122): tscEtoh3 gains a Pathway value
$ws.Range("D121").Value = "ethanol"

# Row 122 (AlleleID 123): tscEtoh4 gains a Pathway value
$ws.Range("D122").Value = "ethanol"

# Row 128 (AlleleID 129): ppdk3 gains a Pathway value
$ws.Range("D128").Value = "malate shunt"

# Row 131 (AlleleID 132): tscEtoh5 -> tscEtoh5int, now an integrated (not
# plasmid-borne) construct
$ws.Range("C131").Value = "tscEtoh5int"
$ws.Range("D131").Value = "ethanol"
$ws.Range("F131").Value = "P2638::pSH62(P2638-adhE(Tsc)-nfnAB(Tsc)-adhA(Tsc))"
$ws.Range("N131").Value = "pSH062 integrated in chromosome at Clo1313_2638 locus"

# Row 132 (AlleleID 133): bcaat1 genotype / note renamed to azlCD naming
$ws.Range("F132").Value = "ΔazlCD"
$ws.Range("N132").Value = "Called in AG1326 and AG2069, also known as ilvT, bcaat, azlCD or brnEF"

# Row 133 (AlleleID 134): bcat1 -> ilvE1 renamed
$ws.Range("C133").Value = "ilvE1"
$ws.Range("F133").Value = "ΔilvE"
$ws.Range("N133").Value = "Called in AG1218, also known as ilvE and bcat"

# ---------------------------------------------------------------------------
# 2. New rows 148-151 (AlleleIDs 150-153)
# ---------------------------------------------------------------------------

# The Start/End columns (I, J) in this sheet are stored as text rather than
# numbers, even though they look numeric. Format those columns as text
# before writing so Excel keeps them as strings, then restore the cell
# style afterwards so no stray number-format style lingers on the cells.
# (Multi-area ranges only apply formatting to the first area here, so set
# each row's I:J pair separately.)
$textCols148 = $ws.Range("I148:J148")
$textCols151 = $ws.Range("I151:J151")
$textCols148.NumberFormat = "@"
$textCols151.NumberFormat = "@"

# Row 148
$ws.Cells.Item(148, 1).Value = 150
$ws.Cells.Item(148, 2).Value = "Clostridium thermocellum"
$ws.Cells.Item(148, 3).Value = "ppdk4"
$ws.Cells.Item(148, 4).Value = "malate shunt"
$ws.Cells.Item(148, 6).Value = "Δppdk::PgapD-cat-hpt"
$ws.Cells.Item(148, 7).Value = $true
$ws.Cells.Item(148, 8).Value = "NC_017304.1"
$ws.Cells.Item(148, 9).Value = "1105704"
$ws.Cells.Item(148, 10).Value = "1107289"
$ws.Cells.Item(148, 11).Value = "Insertion"
$ws.Cells.Item(148, 12).Value = "pZJ03_mero_region"
$ws.Cells.Item(148, 14).Value = "Partial deletion of ppdk using pZJ03, 355 bp removed"
$ws.Cells.Item(148, 21).Value = 1585

# Row 149
$ws.Cells.Item(149, 1).Value = 151
$ws.Cells.Item(149, 2).Value = "Clostridium thermocellum"
$ws.Cells.Item(149, 3).Value = "tscEtoh2"
$ws.Cells.Item(149, 4).Value = "ethanol"
$ws.Cells.Item(149, 6).Value = "pSH64(P2638-nfnAB(Tsc)-adhA(Tsc))"
$ws.Cells.Item(149, 7).Value = $true
$ws.Cells.Item(149, 8).Value = "pSH064"
$ws.Cells.Item(149, 14).Value = "Expression plasmid with T. sacch nfnAB-adhA"

# Row 150
$ws.Cells.Item(150, 1).Value = 152
$ws.Cells.Item(150, 2).Value = "Clostridium thermocellum"
$ws.Cells.Item(150, 3).Value = "tscEtoh5"
$ws.Cells.Item(150, 4).Value = "ethanol"
$ws.Cells.Item(150, 6).Value = "pSH62(P2638-adhE(Tsc)-nfnAB(Tsc)-adhA(Tsc))"
$ws.Cells.Item(150, 7).Value = $true
$ws.Cells.Item(150, 8).Value = "pSH062"
$ws.Cells.Item(150, 14).Value = "Expression plasmid with T. sacch adhE-nfnAB-adhA"

# Row 151
$ws.Cells.Item(151, 1).Value = 153
$ws.Cells.Item(151, 2).Value = "Clostridium thermocellum"
$ws.Cells.Item(151, 3).Value = "bcaat2"
$ws.Cells.Item(151, 4).Value = "amino acid"
$ws.Cells.Item(151, 6).Value = "ΔClo1313_0817-Clo1313_0826"
$ws.Cells.Item(151, 7).Value = $false
$ws.Cells.Item(151, 8).Value = "NC_017304.1"
$ws.Cells.Item(151, 9).Value = "941487"
$ws.Cells.Item(151, 10).Value = "953492"
$ws.Cells.Item(151, 11).Value = "Deletion"
$ws.Cells.Item(151, 13).Value = "Clo1313_0817-Clo1313_0826"
$ws.Cells.Item(151, 14).Value = "Spontaneous deletion of chromosomal region containing bcaat, probably due to recombination of flanking IS200 sequences"
$ws.Cells.Item(151, 21).Value = 12005

# Undo the temporary text number-format so the cells fall back to the
# workbook's default (unstyled) look, matching the rest of the sheet.
$textCols148.Style = "Normal"
$textCols151.Style = "Normal"

# ---------------------------------------------------------------------------
# 3. Update the "Allele" defined name to cover the extended range
# ---------------------------------------------------------------------------
$wb.Names.Item("Allele").RefersTo = "='Allele'!`$A`$1:`$U`$151"
